$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.131.24"
$ws.Range("E2").Value = "  +3.30%  "
$ws.Range("D3").Value = "3.412.50"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.48%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.127"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.14%  "
$ws.Range("E11").Value = "  +4.30%  "
$ws.Range("D12").Value = "3.991.51"
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("E14").Value = "  +5.24%  "
$ws.Range("D15").Value = "3.414.44"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("D17").Value = "62.161.30"
$ws.Range("E17").Value = "  +3.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.94%  "
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.88%  "
$ws.Range("E22").Value = "  +2.27%  "
$ws.Range("D23").Value = "3.545.91"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("E24").Value = "  +15.10%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("E28").Value = "  -5.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +4.39%  "
$ws.Range("E31").Value = "  +4.02%  "
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("D33").Value = "3.444.01"
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "163.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.48%  "
$ws.Range("E42").Value = "  +4.83%  "
$ws.Range("E43").Value = "  +4.01%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.28%  "
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.41%  "
$ws.Range("D50").Value = "2.374.49"
$ws.Range("E50").Value = "  +8.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0265"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.47%  "
